# Draft mapping addition: "Mapping: Spécification métier vers l'extension
# ROR AccomodationFamily" column on the Elements sheet, and a refreshed
# publication Date on the Metadata sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet: bump the "Date" property value.
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# ---------------------------------------------------------------------
# 2. Elements sheet: append a new "Mapping: ..." column (AL) after the
#    existing "Mapping: RIM Mapping" column (AK), matching its look
#    (header style + body style) and fill in the single mapped value.
# ---------------------------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

# Clone formatting (header/body styles, borders, fill...) from the
# neighbouring "Mapping: RIM Mapping" column so the new column matches
# the rest of the mapping columns.
$wsElem.Range("AK1:AK6").Copy()
$wsElem.Range("AL1:AL6").PasteSpecial(-4122)  # xlPasteFormats

# Header for the new mapping column.
$wsElem.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR AccomodationFamily"

# Only the "Extension.value[x]" row (row 6) carries a mapped value; the
# rest of the rows stay blank like the other mapping columns.
$wsElem.Range("AL2:AL5").Value = ""
$wsElem.Range("AL6").Value = "hebergementFamille"

# Widen the new column to fit its (long) header text, like the other
# wide descriptive columns on this sheet.
$wsElem.Columns.Item(38).ColumnWidth = 75.7
